$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newProgramsQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Multiple Myeloma%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

# Update the Programs tab query text in B2 (ProgramsTab row) to add the
# Website CASE expression (program_link / program_acronym).
$ws.Range("B2").Value = $newProgramsQuery

# Leave the active selection on C3, matching where the author's cursor
# ended up after editing the cell.
$ws.Range("C3").Select() | Out-Null
